$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 6")
$ws.Activate()

# Row 14: fill in entry #9
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = 43532
$ws.Range("C14").Value = 0.33333333333333331
$ws.Range("D14").Value = 0.4375
$ws.Range("F14").Formula = "=(D14-C14)*24*60 - E14"
$ws.Range("G14").Value = "Class "
$ws.Range("H14").Value = "Praktikum"

# Row 15: new entry #10
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = 43532
$ws.Range("C15").Value = 0.45833333333333331
$ws.Range("D15").Value = 0.65277777777777779
$ws.Range("E15").Value = 30
$ws.Range("F15").Formula = "=(D15-C15)*24*60 - E15"
$ws.Range("G15").Value = "Proge"
$ws.Range("H15").Value = "MVC EF"

# Row 16: new entry #11
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 43532
$ws.Range("C16").Value = 0.91666666666666663
$ws.Range("D16").Value = 0.97222222222222221
$ws.Range("F16").Formula = "=(D16-C16)*24*60 - E16"
$ws.Range("G16").Value = "Proge"
$ws.Range("H16").Value = "MVC EF"
$ws.Range("I16").Value = "x"

# Row 17: Total row (moved from row 15)
$ws.Range("A17").Value = "Total Time:"
$ws.Range("F17").Formula = "=SUM(F6:F16)"

$ws.Columns.Item(8).ColumnWidth = 13.42578125

$ws.Range("M16").Select()
